$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 2667.8572
$ws.Range("I40").Value = 1294.5
$ws.Range("K40").Value = 1294.5
$ws.Range("M40").Value = -1119.5

# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 12608
$ws.Range("J62").Value = 7000
$ws.Range("L62").Value = 7000
$ws.Range("N62").Value = -8248

# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 12608
$ws.Range("J65").Value = 7000
$ws.Range("L65").Value = 35000
$ws.Range("N65").Value = -41240

# Row 86: Filling in the Blanks
$ws.Range("H86").Value = 2677.5
$ws.Range("I86").Value = 2343
$ws.Range("K86").Value = 2343
$ws.Range("M86").Value = -1220

# Row 89: Ink into Antiquity (L)
$ws.Range("H89").Value = 2677.5
$ws.Range("I89").Value = 2343
$ws.Range("K89").Value = 11715
$ws.Range("M89").Value = -6099

# Row 106: Making Your Mark
$ws.Range("H106").Value = 76526
$ws.Range("I106").Value = 113770
$ws.Range("K106").Value = 113770
$ws.Range("M106").Value = -113139

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 85924.5
$ws.Range("I137").Value = 2399.5
$ws.Range("K137").Value = 7198.5
$ws.Range("M137").Value = -4648.5

# Row 138: All-night Crafting
$ws.Range("H138").Value = 3046.9268
$ws.Range("I138").Value = 1842.875
$ws.Range("J138").Value = 4746.7646
$ws.Range("K138").Value = 5528.625
$ws.Range("L138").Value = 14240.2938
$ws.Range("M138").Value = -388.625
$ws.Range("N138").Value = -24520.2938

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 34368.25
$ws.Range("I141").Value = 34368.25
$ws.Range("K141").Value = 103104.75
$ws.Range("M141").Value = -97924.75

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 147322.89
$ws.Range("I32").Value = 140842.81
$ws.Range("K32").Value = 140842.81
$ws.Range("M32").Value = -140555.81

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 2848.8572
$ws.Range("I61").Value = 2683.3845
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2683.3845
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2471.3845
$ws.Range("N61").Value = -5424

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 2848.8572
$ws.Range("I136").Value = 2683.3845
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 8050.1535
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -5500.1535
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt
$ws.Range("H20").Value = 6180.722
$ws.Range("I20").Value = 6271.6665
$ws.Range("J20").Value = 5998.8335
$ws.Range("K20").Value = 6271.6665
$ws.Range("L20").Value = 5998.8335
$ws.Range("M20").Value = -6024.6665
$ws.Range("N20").Value = -6492.8335

# Row 82: Spirituality Inspector
$ws.Range("H82").Value = 20765.285

# Row 85: The Clamor for Hammers (L)
$ws.Range("H85").Value = 20765.285

# Row 122: To Delight a Dancer
$ws.Range("H122").Value = 78000
$ws.Range("J122").Value = 78000
$ws.Range("L122").Value = 78000
$ws.Range("N122").Value = -87800

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2938.3333
$ws.Range("I134").Value = 2926.2
$ws.Range("K134").Value = 8778.599999999999
$ws.Range("M134").Value = -6243.599999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 4479.049
$ws.Range("I31").Value = 2931.1072
$ws.Range("J31").Value = 7813.077
$ws.Range("K31").Value = 2931.1072
$ws.Range("L31").Value = 7813.077
$ws.Range("M31").Value = -2636.1072
$ws.Range("N31").Value = -8403.077000000001

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 4479.049
$ws.Range("I34").Value = 2931.1072
$ws.Range("J34").Value = 7813.077
$ws.Range("K34").Value = 2931.1072
$ws.Range("L34").Value = 7813.077
$ws.Range("M34").Value = -2729.1072
$ws.Range("N34").Value = -8217.077000000001

# Row 99: O Pine
$ws.Range("H99").Value = 8193.4
$ws.Range("I99").Value = 7991.75
$ws.Range("K99").Value = 7991.75
$ws.Range("M99").Value = -6493.75

# Row 103: Spare a Rod and Spoil the Fishers
$ws.Range("H103").Value = 166679170
$ws.Range("I103").Value = 166679170
$ws.Range("K103").Value = 166679170
$ws.Range("M103").Value = -166677998

# Row 126: A Better Conductor
$ws.Range("H126").Value = 8193.4
$ws.Range("I126").Value = 7991.75
$ws.Range("K126").Value = 23975.25
$ws.Range("M126").Value = -21505.25

$ws = $wb.Worksheets.Item("CUL")
# Row 50: Moving Up in the World
$ws.Range("H50").Value = 143673.42
$ws.Range("J50").Value = 5000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -15962

# Row 53: Rolanberry Fields Forever
$ws.Range("H53").Value = 143673.42
$ws.Range("J53").Value = 5000
$ws.Range("L53").Value = 15000
$ws.Range("N53").Value = -15962

# Row 55: Pagan Pastries
$ws.Range("H55").Value = 4203739
$ws.Range("J55").Value = 5439347.5
$ws.Range("L55").Value = 16318042.5
$ws.Range("N55").Value = -16318396.5

# Row 122: Salt of the North
$ws.Range("H122").Value = 760.2857
$ws.Range("I122").Value = 198
$ws.Range("J122").Value = 854
$ws.Range("K122").Value = 1782
$ws.Range("L122").Value = 7686
$ws.Range("M122").Value = 668
$ws.Range("N122").Value = -12586

# Row 132: More Mezcal
$ws.Range("H132").Value = 2599.6
$ws.Range("J132").Value = 2599.6
$ws.Range("L132").Value = 23396.4
$ws.Range("N132").Value = -28456.4

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 8672.429
$ws.Range("I70").Value = 9699.75
$ws.Range("K70").Value = 9699.75
$ws.Range("M70").Value = -9429.75

# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 8672.429
$ws.Range("I73").Value = 9699.75
$ws.Range("K73").Value = 9699.75
$ws.Range("M73").Value = -8763.75

# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 45462936
$ws.Range("I80").Value = 71438260
$ws.Range("K80").Value = 71438260
$ws.Range("M80").Value = -71437262

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 45462936
$ws.Range("I83").Value = 71438260
$ws.Range("K83").Value = 357191300
$ws.Range("M83").Value = -357186308

# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 33336760
$ws.Range("I113").Value = 38464570
$ws.Range("K113").Value = 38464570
$ws.Range("M113").Value = -38462400

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 19974.834
$ws.Range("I7").Value = 35381.223
$ws.Range("J7").Value = 4568.4443
$ws.Range("K7").Value = 35381.223
$ws.Range("L7").Value = 4568.4443
$ws.Range("M7").Value = -35269.223
$ws.Range("N7").Value = -4792.4443

# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 1162
$ws.Range("I82").Value = 198
$ws.Range("J82").Value = 1483.3334
$ws.Range("K82").Value = 198
$ws.Range("L82").Value = 1483.3334
$ws.Range("M82").Value = 163
$ws.Range("N82").Value = -2205.3334

# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 1162
$ws.Range("I85").Value = 198
$ws.Range("J85").Value = 1483.3334
$ws.Range("K85").Value = 198
$ws.Range("L85").Value = 1483.3334
$ws.Range("M85").Value = 1050
$ws.Range("N85").Value = -3979.3334

# Row 122: Hell on Leather
$ws.Range("H122").Value = 7084.185
$ws.Range("I122").Value = 7490.3477
$ws.Range("K122").Value = 22471.0431
$ws.Range("M122").Value = -20021.0431

# Row 126: Battered Books
$ws.Range("H126").Value = 19974.834
$ws.Range("I126").Value = 35381.223
$ws.Range("J126").Value = 4568.4443
$ws.Range("K126").Value = 106143.669
$ws.Range("L126").Value = 13705.3329
$ws.Range("M126").Value = -103673.669
$ws.Range("N126").Value = -18645.3329

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 3443.24
$ws.Range("I132").Value = 3003.9546
$ws.Range("K132").Value = 9011.863799999999
$ws.Range("M132").Value = -6481.863799999999

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 4509.25
$ws.Range("I136").Value = 3612.0667
$ws.Range("K136").Value = 10836.2001
$ws.Range("M136").Value = -8286.2001

$ws = $wb.Worksheets.Item("WVR")
# Row 3: Trew Enough
$ws.Range("H3").Value = 10048000
$ws.Range("I3").Value = 12545000
$ws.Range("J3").Value = 60000
$ws.Range("K3").Value = 12545000
$ws.Range("L3").Value = 60000
$ws.Range("M3").Value = -12544886
$ws.Range("N3").Value = -60228

# Row 62: Pride Up in Smoke
$ws.Range("H62").Value = 163069.77
$ws.Range("I62").Value = 79979
$ws.Range("J62").Value = 199999
$ws.Range("K62").Value = 79979
$ws.Range("L62").Value = 199999
$ws.Range("M62").Value = -79355
$ws.Range("N62").Value = -201247

# Row 65: Desperate for Diversionaries (L)
$ws.Range("H65").Value = 163069.77
$ws.Range("I65").Value = 79979
$ws.Range("J65").Value = 199999
$ws.Range("K65").Value = 399895
$ws.Range("L65").Value = 999995
$ws.Range("M65").Value = -396775
$ws.Range("N65").Value = -1006235

# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 62501464
$ws.Range("I81").Value = 66668028
$ws.Range("K81").Value = 133336056
$ws.Range("M81").Value = -133334995

# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 62501464
$ws.Range("I84").Value = 66668028
$ws.Range("K84").Value = 666680280
$ws.Range("M84").Value = -666674976

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 7960.9697
$ws.Range("I132").Value = 11237.526
$ws.Range("J132").Value = 3514.2144
$ws.Range("K132").Value = 33712.578
$ws.Range("L132").Value = 10542.6432
$ws.Range("M132").Value = -31182.578
$ws.Range("N132").Value = -15602.6432
